$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AA9" = 29
    "AB9" = 35
    "AC9" = 8.25
    "AD9" = 5.8
    "AG9" = 500
    "AH9" = 7.2
    "AI9" = 11
    "AJ9" = 8.75
    "AK9" = 24
    "AL9" = 19.5
    "AM9" = 29
    "AO9" = 17
    "AP9" = 21
    "AQ9" = 80
    "AR9" = 100
    "AS9" = 200
    "AT9" = 2.6
    "AU9" = 6.3
    "AV9" = 50
    "AW9" = 4.25
    "AX9" = 11.75
    "AY9" = 18
    "G9" = 3.2
    "H9" = 2.95
    "I9" = 2.27
    "J9" = 3.6
    "K9" = 2.05
    "L9" = 2.8
    "M9" = 1.06
    "N9" = 8.5
    "O9" = 1.33
    "P9" = 2.8
    "Q9" = 1.98
    "R9" = 1.65
    "T9" = 2.45
    "V9" = 1.91
    "W9" = 9.5
    "X9" = 17.5
    "Y9" = 11
    "AA27" = 17
    "AH27" = 19
    "AX27" = 12
    "AZ27" = 34
    "G27" = 2.5
    "I27" = 2.5
    "J27" = 2.88
    "N27" = 26
    "O27" = 1.07
    "P27" = 9
    "X27" = 21
    "Z27" = 29
    "AB28" = 19
    "AT28" = 3.75
    "G28" = 1.95
    "H28" = 3.8
    "L28" = 3.6
    "Q28" = 1.48
    "R28" = 2.6
    "S28" = 1.25
    "T28" = 3.75
    "U28" = 1.44
    "V28" = 2.63
    "AC29" = 34
    "AD29" = 15
    "AE29" = 19
    "AG29" = 126
    "AH29" = 29
    "AJ29" = 21
    "AK29" = 81
    "AM29" = 41
    "AO29" = 6
    "AP29" = 12
    "AQ29" = 13
    "AR29" = 26
    "AT29" = 5
    "AU29" = 8
    "AW29" = 9.5
    "AX29" = 34
    "AY29" = 29
    "AZ29" = 101
    "G29" = 1.3
    "H29" = 6.5
    "I29" = 7.5
    "J29" = 1.67
    "K29" = 3
    "L29" = 6.5
    "O29" = 1.07
    "P29" = 9
    "Q29" = 1.29
    "R29" = 3.6
    "S29" = 1.17
    "T29" = 5
    "U29" = 1.53
    "V29" = 2.38
    "Y29" = 9.5
    "Z29" = 11
    "AD30" = 11
    "AH30" = 23
    "AI30" = 34
    "AJ30" = 17
    "AQ30" = 19
    "AR30" = 29
    "AT30" = 4.33
    "G30" = 1.5
    "H30" = 5
    "I30" = 5.25
    "J30" = 1.95
    "K30" = 2.75
    "L30" = 5
    "N30" = 26
    "Q30" = 1.33
    "R30" = 3.4
    "S30" = 1.2
    "T30" = 4.33
    "U30" = 1.44
    "V30" = 2.63
    "W30" = 13
    "Y30" = 9
    "AA31" = 19
    "AC31" = 17
    "AF31" = 29
    "AI31" = 15
    "AJ31" = 10
    "AK31" = 26
    "AL31" = 17
    "AM31" = 21
    "AN31" = 5
    "AO31" = 13
    "AW31" = 4.75
    "AX31" = 13
    "BC31" = 301
    "G31" = 2.7
    "I31" = 2.55
    "J31" = 3.1
    "K31" = 2.4
    "L31" = 3
    "N31" = 15
    "U31" = 1.44
    "V31" = 2.63
    "W31" = 13
    "X31" = 17
    "Y31" = 11
    "Z31" = 29
    "AA32" = 21
    "AB32" = 29
    "AC32" = 11
    "AG32" = 201
    "AH32" = 9.5
    "AI32" = 13
    "AJ32" = 10
    "AL32" = 21
    "AP32" = 23
    "AR32" = 67
    "AS32" = 151
    "AT32" = 2.75
    "AW32" = 4.75
    "G32" = 2.55
    "H32" = 3.4
    "I32" = 2.75
    "J32" = 3.2
    "L32" = 3.25
    "O32" = 1.29
    "P32" = 3.5
    "Q32" = 1.98
    "R32" = 1.88
    "S32" = 1.4
    "T32" = 2.75
    "U32" = 1.7
    "V32" = 2.05
    "X32" = 13
    "Y32" = 10
    "Z32" = 26
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
